{"js": "// Update the header date and the twenty-five two-digit multiplication\n// problems scattered across the table cells. Every \"old\" value is unique\n// in the document, so a simple search-and-replace per pair is safe and\n// unambiguous.\nconst replacements = [\n  [\"2025-11-22 Saturday\", \"2025-11-23 Sunday\"],\n  [\"42\u00d737=\", \"14\u00d761=\"],\n  [\"86\u00d734=\", \"68\u00d750=\"],\n  [\"12\u00d739=\", \"78\u00d719=\"],\n  [\"64\u00d728=\", \"19\u00d793=\"],\n  [\"41\u00d780=\", \"48\u00d734=\"],\n  [\"14\u00d716=\", \"32\u00d784=\"],\n  [\"79\u00d757=\", \"70\u00d792=\"],\n  [\"25\u00d755=\", \"29\u00d748=\"],\n  [\"72\u00d760=\", \"97\u00d771=\"],\n  [\"42\u00d731=\", \"91\u00d793=\"],\n  [\"34\u00d721=\", \"15\u00d717=\"],\n  [\"67\u00d796=\", \"55\u00d736=\"],\n  [\"79\u00d788=\", \"46\u00d722=\"],\n  [\"99\u00d759=\", \"12\u00d792=\"],\n  [\"85\u00d713=\", \"83\u00d766=\"],\n  [\"50\u00d757=\", \"27\u00d786=\"],\n  [\"64\u00d722=\", \"69\u00d787=\"],\n  [\"81\u00d791=\", \"50\u00d796=\"],\n  [\"71\u00d767=\", \"18\u00d762=\"],\n  [\"53\u00d738=\", \"87\u00d757=\"],\n  [\"33\u00d753=\", \"31\u00d758=\"],\n  [\"69\u00d767=\", \"87\u00d774=\"],\n  [\"95\u00d747=\", \"27\u00d744=\"],\n  [\"11\u00d759=\", \"73\u00d720=\"],\n  [\"64\u00d751=\", \"93\u00d712=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the header date and the twenty-five two-digit multiplication\n# problems scattered across the table cells. Every \"old\" value is unique\n# in the document, so a plain Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-11-22 Saturday\", \"2025-11-23 Sunday\"),\n  @(\"42\u00d737=\", \"14\u00d761=\"),\n  @(\"86\u00d734=\", \"68\u00d750=\"),\n  @(\"12\u00d739=\", \"78\u00d719=\"),\n  @(\"64\u00d728=\", \"19\u00d793=\"),\n  @(\"41\u00d780=\", \"48\u00d734=\"),\n  @(\"14\u00d716=\", \"32\u00d784=\"),\n  @(\"79\u00d757=\", \"70\u00d792=\"),\n  @(\"25\u00d755=\", \"29\u00d748=\"),\n  @(\"72\u00d760=\", \"97\u00d771=\"),\n  @(\"42\u00d731=\", \"91\u00d793=\"),\n  @(\"34\u00d721=\", \"15\u00d717=\"),\n  @(\"67\u00d796=\", \"55\u00d736=\"),\n  @(\"79\u00d788=\", \"46\u00d722=\"),\n  @(\"99\u00d759=\", \"12\u00d792=\"),\n  @(\"85\u00d713=\", \"83\u00d766=\"),\n  @(\"50\u00d757=\", \"27\u00d786=\"),\n  @(\"64\u00d722=\", \"69\u00d787=\"),\n  @(\"81\u00d791=\", \"50\u00d796=\"),\n  @(\"71\u00d767=\", \"18\u00d762=\"),\n  @(\"53\u00d738=\", \"87\u00d757=\"),\n  @(\"33\u00d753=\", \"31\u00d758=\"),\n  @(\"69\u00d767=\", \"87\u00d774=\"),\n  @(\"95\u00d747=\", \"27\u00d744=\"),\n  @(\"11\u00d759=\", \"73\u00d720=\"),\n  @(\"64\u00d751=\", \"93\u00d712=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2)\n}\n"}
